$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data in row 3 and row 4 for columns D, J, K, L, M, P
# (the remaining columns already hold identical values in both rows)

$d3 = $ws.Range("D3").Value2
$j3 = $ws.Range("J3").Value2
$k3 = $ws.Range("K3").Value2
$l3 = $ws.Range("L3").Value2
$m3 = $ws.Range("M3").Value2
$p3 = $ws.Range("P3").Value2

$d4 = $ws.Range("D4").Value2
$j4 = $ws.Range("J4").Value2
$k4 = $ws.Range("K4").Value2
$l4 = $ws.Range("L4").Value2
$m4 = $ws.Range("M4").Value2
$p4 = $ws.Range("P4").Value2

$ws.Range("D3").Value = $d4
$ws.Range("J3").Value = $j4
$ws.Range("K3").Value = $k4
$ws.Range("L3").Value = $l4
$ws.Range("M3").Value = $m4
$ws.Range("P3").Value = $p4

$ws.Range("D4").Value = $d3
$ws.Range("J4").Value = $j3
$ws.Range("K4").Value = $k3
$ws.Range("L4").Value = $l3
$ws.Range("M4").Value = $m3
$ws.Range("P4").Value = $p3
